$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace route_ids values that reference the old "purchase" module route
# with the new "purchase_stock" module route (module rename upstream).
for ($r = 1; $r -le $ws.UsedRange.Rows.Count; $r++) {
    $cell = $ws.Cells.Item($r, 18)  # Column R = route_ids
    $val = $cell.Value
    if ($val -ne $null -and $val.ToString().Contains("purchase.route_warehouse0_buy")) {
        $cell.Value = $val.ToString().Replace("purchase.route_warehouse0_buy", "purchase_stock.route_warehouse0_buy")
    }
}

# Update the frozen-pane scroll position / selection state for the sheet view
$ws.Application.ActiveWindow.ScrollRow = 2
$ws.Application.ActiveWindow.ScrollColumn = 18

$ws.Range("A2").Select()
$ws.Range("R2").Activate()

$topLeft = $ws.Range("A1:B1")
$ws.Range("R1:R30").Select()
